# Update the worksheet data: replace rows 2-3 and append rows 4-10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("POST009112-5",   "秀林和平郵局(花蓮12支)",        "花蓮縣秀林鄉和平村113號",              "花蓮"),
    @("POST014122-5",   "和平梨山郵局(台中94支)",        "台中市和平區梨山村中正路89號",          "台中"),
    @("SET180632",      "統一超商  環山店",              "台中市和平區中興路三段64-5、64-6號",    "台中"),
    @("TCAT009203",     "統一速達    梨山衛星所",         "台中市和平區梨山村福壽路5號之1",        "台中"),
    @("TCAT223010",     "統一速達  梨山衛星所",           "台中市和平區梨山村福壽路5號之1",        "台中"),
    @("TCATOA223010",   "統一速達  梨山衛星所",           "台中市和平區梨山村福壽路5號之1",        "台中"),
    @("TRA057",         "台灣鐵路管理局　和平站",         "花蓮縣秀林鄉和平村２７６號",            "花蓮"),
    @("TRA057出入",      "台灣鐵路管理局　和平站  出入口", "花蓮縣秀林鄉和平村２７６號",            "花蓮"),
    @("TRA375",         "台灣鐵路管理局  阿里山站",       "嘉義縣阿里山鄉中正村1號",              "嘉義")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
